$wb = $excel.ActiveWorkbook

# --- "Sheet1" (template sheet): fix/extend the 'header' filter test case ---
# A new template placeholder row is inserted above the existing
# '{{ df2 | noheader | maxrows(...) }}' row so the fixed `header` filter has
# its own dedicated test case, like the other df2 filters on this sheet.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = "{{ df2 | header }}"

# --- "expected" sheet: add the expected rendered output for the new case ---
# The fixed header filter renders df2 with its header row, so the expected
# sheet gets a new header row (matching the df2 header in row 3) right
# before the data rows that follow it.
$ws2 = $wb.Worksheets.Item("expected")
$ws2.Rows.Item(7).Insert()
$ws2.Range("A7").Value = "index"
$ws2.Range("B7").Value = "name"
$ws2.Range("C7").Value = "b"
$ws2.Range("D7").Value = "c"
$ws2.Range("E7").Value = "d"

# --- restore selections / active sheet as left by the author ---
$ws2.Activate()
$ws2.Range("B17").Select()
$ws1.Activate()
$ws1.Range("A15").Select()
